$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 01:05"

# Rows whose ranking shuffled: a new/fast-growing country bumped ahead of its
# neighbours, so the country name (column A) for these rows changes along
# with the statistics (B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes).
$rowNames = @{
    35 = "Colombia"
    36 = "Indonesia"
    37 = "Kuwait"
    49 = "Panama"
    50 = "Dinamarca"
    57 = "Nigeria"
    58 = "Noruega"
    59 = "Oman"
    64 = "Ghana"
    65 = "Australia"
    66 = "Bolivia"
}

foreach ($rowNum in $rowNames.Keys) {
    $ws.Range("A$rowNum").Value = $rowNames[$rowNum]
}

# Updated per-row statistics
$updates = @{
    4  = @{ B = 1744761; C = 19504; D = 486253; E = 1156461; G = 1475; H = 102047 }
    16 = @{ B = 87519;   C = 872;   D = 46164;  E = 34590;   G = 126;  H = 6765 }
    35 = @{ B = 24104;   C = 1101;  D = 6111;   E = 17190;   G = 27;   H = 803 }
    36 = @{ B = 23851;   C = 686;   D = 6057;   E = 16321;   G = 55;   H = 1473 }
    37 = @{ B = 23267;   C = 692;   D = 7946;   E = 15146;   G = 3;    H = 175 }
    49 = @{ B = 11728;   C = 281;   D = 7379;   E = 4034;    H = 315 }
    50 = @{ B = 11480;   C = 52;    D = 10106;  E = 809;     G = 2;    H = 565 }
    55 = @{ B = 9086;    C = 36;    D = 6370;   E = 2399 }
    57 = @{ B = 8733;    C = 389;   D = 2501;   E = 5978;    G = 5;    H = 254 }
    58 = @{ B = 8401;    C = 18;    D = 7727;   E = 439;     G = 0;    H = 235 }
    59 = @{ B = 8373;    C = 255;   D = 2177;   E = 6157;    G = 2;    H = 39 }
    64 = @{ B = 7303;    C = 186;   D = 2412;   E = 4857;    G = 0;    H = 34 }
    65 = @{ B = 7139;    C = 6;     D = 6566;   E = 470;     G = 1;    H = 103 }
    66 = @{ B = 7136;    C = 476;   D = 677;    E = 6185;    G = 13;   H = 274 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
